$d = $word.ActiveDocument

$replacements = @(
    @("682×3=2046", "798×3=2394"),
    @("282×3=846", "336×3=1008"),
    @("911×5=4555", "104×2=208"),
    @("767×8=6136", "607×4=2428"),
    @("762×7=5334", "566×3=1698"),
    @("828×4=3312", "606×5=3030"),
    @("631×2=1262", "374×3=1122"),
    @("354×6=2124", "170×7=1190"),
    @("657×9=5913", "840×4=3360"),
    @("508×7=3556", "431×7=3017"),
    @("177×9=1593", "402×3=1206"),
    @("902×4=3608", "455×9=4095"),
    @("241×4=964", "771×4=3084"),
    @("442×9=3978", "911×7=6377"),
    @("238×3=714", "496×4=1984"),
    @("743×8=5944", "655×8=5240"),
    @("508×6=3048", "724×8=5792"),
    @("999×6=5994", "538×2=1076"),
    @("256×6=1536", "262×9=2358"),
    @("384×3=1152", "449×7=3143"),
    @("556×2=1112", "708×5=3540"),
    @("569×7=3983", "691×5=3455"),
    @("693×7=4851", "912×9=8208"),
    @("877×7=6139", "669×3=2007"),
    @("465×8=3720", "879×4=3516")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
